# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 163 (01-06-2021) with revised figures
$ws.Range("B163").Value = 980835
$ws.Range("C163").Value = 10191
$ws.Range("D163").Value = 805553
$ws.Range("E163").Value = 115488
$ws.Range("F163").Value = 59794

# Add new row 164 for 01-07-2021
# Enter the date-like label as a formula-literal then collapse it to a
# plain value via copy/paste-values, so it lands as shared-string text
# instead of being auto-converted into a date serial number/format.
$ws.Range("A164").Formula = '="01-07-2021"'
$ws.Range("A164").Copy()
$ws.Range("A164").PasteSpecial(-4163)
$ws.Range("B164").Value = 1317319
$ws.Range("C164").Value = 12196
$ws.Range("D164").Value = 1033024
$ws.Range("E164").Value = 223527
$ws.Range("F164").Value = 60768
